$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the row for account 004940699 / RACHEL / 17138.49 (row 10)
$ws.Rows.Item(10).Delete()

# 2) Insert a new row (004958578 / ASSAKO / 26.36) right before the
#    001719494 / LUIS row, which is now row 155 after the deletion above.
$ws.Rows.Item(155).Insert()
$ws.Cells.Item(155, 1).NumberFormat = "@"
$ws.Cells.Item(155, 1).Value = "004958578"
$ws.Cells.Item(155, 2).Value = "ASSAKO"
$ws.Cells.Item(155, 3).Value = 26.36

# 3) Remove the trailing four rows (BRUNO, ASSAKO -725.91, PATRICIA, CEZAR)
#    that used to sit right after the VILMA 0.01 row (rows 331-334, still
#    at the same indices since edits above were both above this block).
$ws.Range("A331:C334").EntireRow.Delete()
